# Weekly update: insert a new price row (week of 2021-09-28, serial 44467)
# at row 51, pushing all existing rows 51-114 down to 52-115.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 51 (shifts rows 51:114 -> 52:115,
# and the used range/dimension grows from R114 to R115 automatically).
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new weekly record.
$ws.Range("A51").Value = 11
$ws.Range("B51").Value = 'Vega Monumental Concepción'
$ws.Range("C51").Value = 'Bíobío'
$ws.Range("D51").Value2 = 44467
$ws.Range("E51").Value = 8
$ws.Range("F51").Value = 100114001
$ws.Range("G51").Value = 'Papa'
$ws.Range("H51").Value = 'Asterix'
$ws.Range("I51").Value = '1a (guarda)'
$ws.Range("J51").Value = 2000
$ws.Range("K51").Value = 9500
$ws.Range("L51").Value = 10000
$ws.Range("M51").Value = 9750
$ws.Range("N51").Value = '$/saco 25 kilos'
$ws.Range("O51").Value = 'Provincia de Arauco'
$ws.Range("P51").Value = 390
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = 'Hortaliza'
